# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the first data row
# (the 7255c0d1... entry) on both the zh-cn and de-de worksheets, to
# reflect a newer handback run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-18 20:40:49"
$zhcn.Range("H2").Value = "2016-03-18 20:41:08"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-18 20:40:52"
$dede.Range("H2").Value = "2016-03-18 20:41:13"
